# Applies the edit described by the diff:
#   "... the first 3 common venue doesn't include restaurants ..."
# becomes
#   "... the first 3 common venue include restaurants ..."
# i.e. the word "doesn't " is deleted from the bold "clusters ..." sentence,
# and Word's automatic "_GoBack" last-edit bookmark moves from its old spot
# (right after "Conclusion:" at the end of the doc) to the new edit point
# (right after "venue ").

$d = $word.ActiveDocument

$apos = [char]0x2019   # curly apostrophe used in "doesn't"

# Locate "doesn't " (with trailing space) in the document.
$target = $d.Content
$target.Find.Execute("doesn" + $apos + "t ", $true, $false, $false, $false, $false, $true, 1, $false, "", $false)

if ($target.Find.Found) {
    $delStart = $target.Start
    $delEnd = $target.End

    # Drop a temporary bookmark right before "doesn't " so the runtime
    # doesn't coalesce this run with the preceding one while we edit -
    # this preserves the existing run split exactly as Word would.
    $barrier = $d.Range($delStart, $delStart)
    $d.Bookmarks.Add("zzTempBarrier", $barrier)

    # "doesn't " splits as "doesn'" (6 chars) + "t " (2 chars); the real
    # edit point (and new _GoBack position) sits right between them.
    $bmPos = $delStart + 6
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    # Remove "t " immediately after the bookmark.
    $afterBm = $d.Range($bmPos, $bmPos + 2)
    $afterBm.Text = ""

    # Remove "doesn'" immediately before the bookmark.
    $beforeBm = $d.Range($delStart, $delStart + 6)
    $beforeBm.Text = ""

    # Clean up the temporary barrier bookmark.
    $d.Bookmarks.Item("zzTempBarrier").Delete()
}
